# Fill in the two remaining placeholder password cells on the "Login Details"
# sheet with the actual test-user's surname/first name (Cele / Nkosi), as
# part of wiring up the Excel-backed test data for the Selenium suite.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login Details")

# B5 was "secret_sauce" -> now holds the surname "Cele"
$ws.Range("B5").Value = "Cele"

# B6 is left untouched (still "secret_sauce")

# B7 was "secret_sauce" -> now holds the first name "Nkosi"
$ws.Range("B7").Value = "Nkosi"

# Reflect the author's on-screen state when they saved the file: the
# "Login Details" sheet was scrolled/zoomed in and the selection left on
# a cell further down the (otherwise empty) sheet.
$ws.Activate()
$excel.ActiveWindow.Zoom = 180
$ws.Range("E612").Select()
